$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated C/D column values (td_sim_1 / record_atd) for corrected relevance markers
$ws.Range("C3").Value = 314
$ws.Range("D3").Value = 314
$ws.Range("C5").Value = 101
$ws.Range("D5").Value = 101
$ws.Range("C7").Value = 65
$ws.Range("D7").Value = 65
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = 32
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 10
$ws.Range("C13").Value = 90
$ws.Range("D13").Value = 90
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 30
$ws.Range("C17").Value = 72
$ws.Range("D17").Value = 72
$ws.Range("C19").Value = 129
$ws.Range("D19").Value = 129
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 5
$ws.Range("C23").Value = 305
$ws.Range("D23").Value = 305
$ws.Range("C25").Value = 85
$ws.Range("D25").Value = 85
$ws.Range("C27").Value = 204
$ws.Range("D27").Value = 204
$ws.Range("C29").Value = 119
$ws.Range("D29").Value = 119
$ws.Range("C31").Value = 385
$ws.Range("D31").Value = 385
$ws.Range("C33").Value = 11
$ws.Range("D33").Value = 11
$ws.Range("C36").Value = 19
$ws.Range("D36").Value = 19
$ws.Range("C38").Value = 89
$ws.Range("D38").Value = 89
$ws.Range("C40").Value = 103
$ws.Range("D40").Value = 103
$ws.Range("C43").Value = 33
$ws.Range("D43").Value = 33
$ws.Range("C45").Value = 20
$ws.Range("D45").Value = 20
$ws.Range("C47").Value = 188
$ws.Range("D47").Value = 188
$ws.Range("C51").Value = 73
$ws.Range("D51").Value = 73
$ws.Range("C53").Value = 1056
$ws.Range("D53").Value = 1056
$ws.Range("C55").Value = 18
$ws.Range("D55").Value = 18
$ws.Range("C57").Value = 35
$ws.Range("D57").Value = 35
$ws.Range("C59").Value = 123
$ws.Range("D59").Value = 123
$ws.Range("C61").Value = 25
$ws.Range("D61").Value = 25
$ws.Range("C63").Value = 22
$ws.Range("D63").Value = 22
$ws.Range("C65").Value = 1
$ws.Range("D65").Value = 1

# Updated average of td_sim_1 (C66) to reflect corrected values
$ws.Range("C66").Value = 121.3030303030303
